$d = $word.ActiveDocument

# 1. Remove the trailing space at the end of the "Segundo teorema" paragraph.
$d.Content.Find.Execute("agir. ", $true, $false, $false, $false, $false, $true, 1, $false, "agir.", 2) | Out-Null

# 2. Append the new "Preco-Sombra" / "Indice de Lerner" content as new paragraphs
#    right before the final paragraph mark (so the existing last paragraph's
#    text/formatting is preserved, and no stray runs are introduced).
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newParagraphsXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Preço-Sombra</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  No planejamento governamental: Se um governo está alocando recursos limitados para diferentes projetos (como infraestrutura, saúde, educação), o preço-sombra ajuda a determinar qual seria o custo de desviar recursos de um projeto para outro.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  Em empresas: Se uma empresa tem uma máquina que pode ser utilizada para diferentes produtos, o preço-sombra indicaria o valor adicional de produzir um produto em detrimento de outro, caso a capacidade de produção seja limitada</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Índice de Lerner (Poder de monopólio)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L = (P – Cmg)/ P</w:t></w:r><w:r><w:t xml:space="preserve">   ou </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>L = -1/E</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:vertAlign w:val="subscript"/></w:rPr><w:t>p</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  E</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="subscript"/></w:rPr><w:t>p</w:t></w:r><w:r><w:t xml:space="preserve"> : Elasticidade preço da demanda</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">  0 &lt;= L &lt;= 1</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($newParagraphsXml)
